$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C49").Value = "'45"
$ws.Range("D49").Value = "'122753.26"
$ws.Range("C50").Value = "'110"
$ws.Range("D50").Value = "'282968.33"
$ws.Range("C51").Value = "'17"
$ws.Range("D51").Value = "'56100.00"
$ws.Range("C52").Value = "'253"
$ws.Range("D52").Value = "'910903.67"
$ws.Range("C55").Value = "'5"
$ws.Range("D55").Value = "'11500.00"
$ws.Range("C56").Value = "'37"
$ws.Range("D56").Value = "'96000.00"
$ws.Range("C57").Value = "'14"
$ws.Range("D57").Value = "'39000.00"
$ws.Range("C60").Value = "'39"
$ws.Range("D60").Value = "'158456.00"
$ws.Range("C61").Value = "'65"
$ws.Range("D61").Value = "'140500.00"
$ws.Range("C74").Value = "'14"
$ws.Range("D74").Value = "'55000.00"
$ws.Range("C77").Value = "'81"
$ws.Range("D77").Value = "'214487.00"
$ws.Range("C79").Value = "'23"
$ws.Range("D79").Value = "'67991.00"
$ws.Range("C80").Value = "'445"
$ws.Range("D80").Value = "'1831304.99"
$ws.Range("C86").Value = "'42"
$ws.Range("D86").Value = "'97500.00"
$ws.Range("C89").Value = "'100"
$ws.Range("D89").Value = "'250197.00"
$ws.Range("C105").Value = "'9"
$ws.Range("D105").Value = "'25909.00"
$ws.Range("C106").Value = "'17"
$ws.Range("D106").Value = "'47709.84"
$ws.Range("C107").Value = "'58"
$ws.Range("D107").Value = "'145860.00"
$ws.Range("C108").Value = "'28"
$ws.Range("D108").Value = "'83934.00"
$ws.Range("C109").Value = "'12"
$ws.Range("D109").Value = "'49913.61"
$ws.Range("C110").Value = "'69"
$ws.Range("D110").Value = "'420236.29"
$ws.Range("C113").Value = "'20"
$ws.Range("D113").Value = "'60267.00"
$ws.Range("C114").Value = "'22"
$ws.Range("D114").Value = "'55895.00"
$ws.Range("C115").Value = "'12"
$ws.Range("D115").Value = "'29100.00"
$ws.Range("C117").Value = "'16"
$ws.Range("D117").Value = "'81445.92"
$ws.Range("C118").Value = "'23"
$ws.Range("D118").Value = "'77197.00"
$ws.Range("C120").Value = "'30"
$ws.Range("D120").Value = "'106000.00"
$ws.Range("C121").Value = "'62"
$ws.Range("D121").Value = "'170877.00"
$ws.Range("C122").Value = "'232"
$ws.Range("D122").Value = "'633208.00"
$ws.Range("C123").Value = "'86"
$ws.Range("D123").Value = "'236571.45"
$ws.Range("C124").Value = "'443"
$ws.Range("D124").Value = "'1888880.06"
$ws.Range("C128").Value = "'83"
$ws.Range("D128").Value = "'246743.68"
$ws.Range("C131").Value = "'19"
$ws.Range("D131").Value = "'39500.00"
$ws.Range("C132").Value = "'77"
$ws.Range("D132").Value = "'335586.75"
$ws.Range("C133").Value = "'113"
$ws.Range("D133").Value = "'282626.44"
$ws.Range("C138").Value = "'552"
$ws.Range("D138").Value = "'1371046.00"
$ws.Range("C139").Value = "'1752"
$ws.Range("D139").Value = "'4685039.93"
$ws.Range("C140").Value = "'2271"
$ws.Range("D140").Value = "'5666808.29"
$ws.Range("C141").Value = "'2424"
$ws.Range("D141").Value = "'10078623.38"
$ws.Range("C142").Value = "'342"
$ws.Range("D142").Value = "'964007.04"
$ws.Range("C145").Value = "'996"
$ws.Range("D145").Value = "'2591433.25"
$ws.Range("C146").Value = "'466"
$ws.Range("D146").Value = "'1345454.49"
$ws.Range("C149").Value = "'385"
$ws.Range("D149").Value = "'1196090.28"
$ws.Range("C150").Value = "'827"
$ws.Range("D150").Value = "'1982242.82"
$ws.Range("C191").Value = "'51"
$ws.Range("D191").Value = "'154300.00"
$ws.Range("C194").Value = "'346"
$ws.Range("D194").Value = "'936288.00"
$ws.Range("C195").Value = "'31"
$ws.Range("D195").Value = "'101574.12"
$ws.Range("C196").Value = "'616"
$ws.Range("D196").Value = "'2253053.66"
$ws.Range("C199").Value = "'31"
$ws.Range("D199").Value = "'77000.00"
$ws.Range("C201").Value = "'60"
$ws.Range("D201").Value = "'169926.00"
$ws.Range("C202").Value = "'75"
$ws.Range("D202").Value = "'177005.00"
$ws.Range("C204").Value = "'113"
$ws.Range("D204").Value = "'512180.50"
$ws.Range("C205").Value = "'132"
$ws.Range("D205").Value = "'292696.77"
$ws.Range("C236").Value = "'21"
$ws.Range("D236").Value = "'63250.00"
$ws.Range("C237").Value = "'78"
$ws.Range("D237").Value = "'212538.00"
$ws.Range("C238").Value = "'142"
$ws.Range("D238").Value = "'368200.00"
$ws.Range("C239").Value = "'484"
$ws.Range("D239").Value = "'1250575.83"
$ws.Range("C240").Value = "'82"
$ws.Range("D240").Value = "'240627.11"
$ws.Range("C241").Value = "'937"
$ws.Range("D241").Value = "'3329351.77"
$ws.Range("C242").Value = "'38"
$ws.Range("D242").Value = "'99500.00"
$ws.Range("C243").Value = "'24"
$ws.Range("D243").Value = "'60000.00"
$ws.Range("C244").Value = "'77"
$ws.Range("D244").Value = "'176500.00"
$ws.Range("C245").Value = "'179"
$ws.Range("D245").Value = "'558429.19"
$ws.Range("C246").Value = "'120"
$ws.Range("D246").Value = "'381293.00"
$ws.Range("C247").Value = "'89"
$ws.Range("D247").Value = "'237972.92"
$ws.Range("C248").Value = "'29"
$ws.Range("D248").Value = "'68500.00"
$ws.Range("C249").Value = "'116"
$ws.Range("D249").Value = "'407662.14"
$ws.Range("C250").Value = "'204"
$ws.Range("D250").Value = "'447013.00"